# "MAS Arreglada para igual"
# Row 8 (the "FinSimboloIncPos" action row) gets a dedicated action name,
# "FinSimbolo", everywhere except the I8 cell (which legitimately keeps the
# longer "FinSimboloIncPos" action). The lone odd-one-out cell O8 also had a
# redundant/duplicate cell style (fill explicitly re-applied with no actual
# fill color) that differed from its row siblings only in bookkeeping, not
# in appearance - clear that stray fill so it matches the rest of the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B8:H8 and J8:R8 -> "FinSimbolo" (I8 stays "FinSimboloIncPos")
$ws.Range("B8:H8").Value = "FinSimbolo"
$ws.Range("J8:R8").Value = "FinSimbolo"

# O8 had a stray "apply fill" flag (no visible fill) that none of the other
# cells on the row carry - drop it so the cell's formatting lines up with
# the rest of row 8.
$ws.Range("O8").Interior.Pattern = -4142

# Zoom out a bit on the sheet (62%).
$excel.ActiveWindow.Zoom = 62

# Leave the selection on the row that was just edited.
$ws.Range("B8:R8").Select()
